$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.294508934020996
$ws.Range("B1").Value = 2.327085971832275
$ws.Range("C1").Value = 2.987038850784302
$ws.Range("D1").Value = 3.439496040344238
$ws.Range("E1").Value = 1.472819447517395
